$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.617.99"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "1.596.42"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'211.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").Value = "'0.515"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.05%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("D10").Value = "'19.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.25%  "
$ws.Range("D11").Value = "'0.0838"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.43%  "
$ws.Range("D12").Value = "1.820.27"
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.578.08"
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'4.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.05%  "
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").Value = "'64.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("D17").Value = "26.607.43"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "0.0₃0732"
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("D19").Value = "'208.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").Value = "'6.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.52%  "
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("D23").Value = "'2.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.00%  "
$ws.Range("D24").Value = "'8.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.34%  "
$ws.Range("D25").Value = "'145.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.08%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "'7.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.49%  "
$ws.Range("D28").Value = "'0.114"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.36%  "
$ws.Range("D29").Value = "'15.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("E32").Value = "  -0.12%  "
$ws.Range("D33").Value = "'0.655"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.64%  "
$ws.Range("D34").Value = "'2.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.26%  "
$ws.Range("D35").Value = "1.279.55"
$ws.Range("E35").Value = "  -2.03%  "
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("E37").Value = "  +0.94%  "
$ws.Range("E38").Value = "  -0.39%  "
$ws.Range("E39").Value = "  +1.72%  "
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("E41").Value = "  +1.95%  "
$ws.Range("D42").Value = "'64.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.82%  "
$ws.Range("E43").Value = "  -0.37%  "
$ws.Range("E44").Value = "  +1.23%  "
$ws.Range("D45").Value = "1.732.91"
$ws.Range("E45").Value = "  +0.56%  "
$ws.Range("E46").Value = "  +8.45%  "
$ws.Range("D47").Value = "'89.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.21%  "
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "'0.103"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.01%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0507"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.46%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.49"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.63%  "
